$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 5
$ws.Range("F2").Value = 5
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 5
$ws.Range("E3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 5
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 5
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("C5").Value = 4
$ws.Range("E5").Value = 4
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 3
$ws.Range("K5").Value = 4
$ws.Range("O5").Value = 4
$ws.Range("C6").Value = 2
$ws.Range("E6").Value = 3
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 5
$ws.Range("K6").Value = 2
$ws.Range("O6").Value = 2
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 3
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 2
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 3
$ws.Range("N7").Value = 3
$ws.Range("O7").Value = 3
